# --array-split option extended to N dimensions
# Splits the previously single-line "a,b" style sample values in H2/G3
# (and the former plain "H3") into multi-line ("\n"-separated) values so
# the sample sheet demonstrates the new N-dimensional array-split syntax.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$nl = [char]10

# ---------------------------------------------------------------------
# 1. Update the sample values to show multi-dimensional splitting
# ---------------------------------------------------------------------
$ws.Range("H2").Value = "H2a1,H2b1" + $nl + "H2a2,H2b2"
$ws.Range("G3").Value = "G3a1,G3b1" + $nl + "G3a2"
$ws.Range("H3").Value = "H3a1" + $nl + "H3a2"

# ---------------------------------------------------------------------
# 2. Give every populated cell a centered vertical alignment, and turn
#    wrap-text on for the new multi-line cells so the embedded line
#    breaks actually render as separate lines.
# ---------------------------------------------------------------------
$vCenterCells = @("A1","B2","G2","B3","B4","B5","H5","A6","B7","G7","B8","G8","A9","B10","G10","B11","G11")
foreach ($addr in $vCenterCells) {
    $ws.Range($addr).VerticalAlignment = -4108
}

$wrapCells = @("H2","G3","H3")
foreach ($addr in $wrapCells) {
    $ws.Range($addr).VerticalAlignment = -4108
    $ws.Range($addr).WrapText = $true
}

# ---------------------------------------------------------------------
# 3. Widen columns G:H slightly and bump rows 2:3 so the two-line values
#    are fully visible, and tighten A:F a touch to match the refreshed
#    layout.
# ---------------------------------------------------------------------
$ws.Columns("A:F").ColumnWidth = 2.0
$ws.Columns("G:H").ColumnWidth = 12.785714285714286

$ws.Rows("2:3").RowHeight = 36

# ---------------------------------------------------------------------
# 4. Move the active selection to G4, matching where the cursor ended up
# ---------------------------------------------------------------------
$ws.Range("G4").Select()
